$d = $word.ActiveDocument

# Helper: create a zero-length bookmark at an absolute document character
# offset. Word (and this COM-interop host) treats a bookmark boundary as a
# place where adjacent, identically-formatted runs are *not* coalesced back
# into one run, which lets us reproduce the exact run layout the diff
# expects even though plain text edits get their containing paragraph
# re-flattened to the minimal run set. The bookmark itself is removed
# immediately after, leaving no trace in the saved document.
function New-SplitMark([int]$pos, [string]$name) {
    $r = $d.Range($pos, $pos)
    $d.Bookmarks.Add($name, $r) | Out-Null
}

# ------------------------------------------------------------------
# Edit 1: in the "Repelling particle ..." paragraph only, change
#   "We do not usually depict those kind of particles ..."
# to
#   "We do not usually depict this kind of particles ..."
# (the identical sentence that also appears in the earlier
# "Match-seeker particle ..." paragraph must stay untouched, so we
# locate this specific paragraph first and work within its range).
# ------------------------------------------------------------------
$repelPara = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -match "Repelling particle is denoted with R\.") {
        $repelPara = $p
        break
    }
}

if ($repelPara -ne $null) {
    $pStart = $repelPara.Range.Start
    $pText = $repelPara.Range.Text
    $thoseIdx = $pText.IndexOf("those")
    if ($thoseIdx -ge 0) {
        # "those" = "th" + "ose" -> keep "th", replace "ose" with "is",
        # matching the run boundary right after "...depict th".
        $oseStart = $pStart + $thoseIdx + 2
        $oseEnd = $oseStart + 3

        New-SplitMark $oseStart "zzSplitA"
        New-SplitMark $oseEnd "zzSplitB"

        $rng = $d.Range($oseStart, $oseEnd)
        $rng.Text = "is"

        $d.Bookmarks("zzSplitA").Delete()
        $d.Bookmarks("zzSplitB").Delete()
    }
}

# ------------------------------------------------------------------
# Edit 2: in the "To each `key-match` property of V ..." paragraph
# only, change
#   "` property of V "
# to
#   "` property of object particle V "
# by inserting "object particle " right before the "V ", splitting the
# original single run into three: "` property of ", "object particle ",
# and "V " -- exactly as in the diff. (Two other, differently-worded
# paragraphs elsewhere also mention `key-match`; matching on the
# trailing "property of V" keeps this scoped to the right one.)
# ------------------------------------------------------------------
$keyMatchPara = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -match "property of V ") {
        $keyMatchPara = $p
        break
    }
}

if ($keyMatchPara -ne $null) {
    $pStart2 = $keyMatchPara.Range.Start
    $pText2 = $keyMatchPara.Range.Text
    $propIdx = $pText2.IndexOf("property of V")
    if ($propIdx -ge 0) {
        $vIdx = $pText2.IndexOf("V", $propIdx)
        $vStart = $pStart2 + $vIdx

        # Mark the insertion point so "` property of " (before) does not
        # re-absorb the text we are about to insert ...
        New-SplitMark $vStart "zzSplitC"

        $insRng = $d.Range($vStart, $vStart)
        $insRng.InsertBefore("object particle ")

        # ... and mark the (now shifted) start of the inserted text so it
        # does not re-absorb into "` property of " on the left either.
        New-SplitMark $vStart "zzSplitD"

        $d.Bookmarks("zzSplitC").Delete()
        $d.Bookmarks("zzSplitD").Delete()
    }
}
